$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "74.569.00"
$ws.Range("E2").Value = "  +6.44%  "
$ws.Range("D3").Value = "2.654.51"
$ws.Range("E3").Value = "  +8.03%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "185.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +10.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "583.32"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.65%  "
$ws.Range("E8").Value = "  +3.60%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.192"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +8.94%  "
$ws.Range("D10").Value = "2.654.47"
$ws.Range("E10").Value = "  +8.16%  "
$ws.Range("E11").Value = "  +1.12%  "
$ws.Range("E12").Value = "  +5.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.73"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.63%  "
$ws.Range("D14").Value = "3.144.07"
$ws.Range("E14").Value = "  +7.97%  "
$ws.Range("D15").Value = "74.416.21"
$ws.Range("E15").Value = "  +6.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000185"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.42%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +8.52%  "
$ws.Range("D18").Value = "2.658.58"
$ws.Range("E18").Value = "  +7.95%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +29.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +9.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "371.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +8.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +10.73%  "
$ws.Range("E23").Value = "  +4.14%  "
$ws.Range("E24").Value = "  +3.10%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "69.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.57%  "
$ws.Range("D29").Value = "2.794.29"
$ws.Range("E29").Value = "  +7.80%  "
$ws.Range("E30").Value = "  -8.42%  "
$ws.Range("D31").Value = "0.0₃0932"
$ws.Range("E31").Value = "  +7.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.41"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +12.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "517.49"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +12.85%  "
$ws.Range("E34").Value = "  +2.94%  "
$ws.Range("E35").Value = "  +6.31%  "
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "162.58"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.118"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.16"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.35"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.36%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "168.84"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +25.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.94"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +11.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.326"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.66"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.17%  "
$ws.Range("E46").Value = "  +7.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "38.97"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.66%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.33"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.52%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0841"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +15.60%  "
$ws.Range("E50").Value = "  +5.88%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.11"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +20.69%  "
